# Updated cryptos list - applies price/volume/link changes per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.398.44'
$ws.Range("E2").Value = '  +0.83%  '

$ws.Range("D3").Value = '2.193.74'
$ws.Range("E3").Value = '  +0.28%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '''254.38'
$ws.Range("E5").Value = '  +5.54%  '

$ws.Range("E6").Value = '  +1.19%  '

$ws.Range("D7").Value = '''68.68'
$ws.Range("E7").Value = '  -2.00%  '

$ws.Range("E8").Value = '  +0.14%  '

$ws.Range("D9").Value = '''0.587'
$ws.Range("E9").Value = '  +9.10%  '

$ws.Range("D10").Value = '''38.09'
$ws.Range("E10").Value = '  +5.09%  '

$ws.Range("D11").Value = '''58.78'
$ws.Range("E11").Value = '  +1.93%  '

$ws.Range("D12").Value = '''0.0941'
$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("D13").Value = '''7.17'
$ws.Range("E13").Value = '  +9.47%  '

$ws.Range("E14").Value = '  +0.65%  '

$ws.Range("D15").Value = '2.516.69'
$ws.Range("E15").Value = '  +0.46%  '

$ws.Range("D16").Value = '''0.875'
$ws.Range("E16").Value = '  +5.20%  '

$ws.Range("D17").Value = '''14.61'
$ws.Range("E17").Value = '  +0.40%  '

$ws.Range("D18").Value = '2.195.12'
$ws.Range("E18").Value = '  +0.94%  '

$ws.Range("D19").Value = '41.262.49'
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").Value = '0.0₃0957'
$ws.Range("E20").Value = '  +2.27%  '

$ws.Range("D21").Value = '''6.25'
$ws.Range("E21").Value = '  +3.95%  '

$ws.Range("D22").Value = '''72.25'
$ws.Range("E22").Value = '  -1.71%  '

$ws.Range("D23").Value = '''233.26'
$ws.Range("E23").Value = '  +1.22%  '

$ws.Range("D24").Value = '''2.06'
$ws.Range("E24").Value = '  +2.44%  '

$ws.Range("D25").Value = '''11.94'
$ws.Range("E25").Value = '  +23.54%  '

$ws.Range("E26").Value = '  +8.44%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("D28").Value = '''2.52'
$ws.Range("E28").Value = '  +5.32%  '

$ws.Range("E29").Value = '  -0.50%  '

$ws.Range("D30").Value = '''170.03'
$ws.Range("E30").Value = '  +1.25%  '

$ws.Range("D31").Value = '''20.72'
$ws.Range("E31").Value = '  +2.80%  '

$ws.Range("E32").Value = '  +2.47%  '

$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").Value = '''5.56'
$ws.Range("E33").Value = '  +9.61%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").Value = '''0.124'
$ws.Range("E34").Value = '  +0.62%  '

$ws.Range("D35").Value = '''0.0738'
$ws.Range("E35").Value = '  +5.88%  '

$ws.Range("D36").Value = '''26.72'
$ws.Range("E36").Value = '  +14.53%  '

$ws.Range("D37").Value = '''4.64'
$ws.Range("E37").Value = '  +2.04%  '

$ws.Range("D38").Value = '''4.11'
$ws.Range("E38").Value = '  +7.15%  '

$ws.Range("E39").Value = '  +13.99%  '

$ws.Range("D40").Value = '''12.85'
$ws.Range("E40").Value = '  +28.34%  '

$ws.Range("E41").Value = '  -0.75%  '

$ws.Range("E42").Value = '  -0.46%  '

$ws.Range("D43").Value = '''64.57'
$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("D44").Value = '''5.01'
$ws.Range("E44").Value = '  +3.68%  '

$ws.Range("D45").Value = '''0.204'
$ws.Range("E45").Value = '  +7.14%  '

$ws.Range("D46").Value = '''8.67'
$ws.Range("E46").Value = '  +0.98%  '

$ws.Range("E47").Value = '  +3.57%  '

$ws.Range("E48").Value = '  +0.19%  '

$ws.Range("E49").Value = '  +5.53%  '

$ws.Range("E50").Value = '  +2.38%  '

$ws.Range("D51").Value = '''4.33'
$ws.Range("E51").Value = '  -4.00%  '

